# Workbook "data/users.xlsx": fill in the password column (C) for each
# user row with sequential numeric codes, except for row 22 (Alireza
# Sheikh al-Eslami) who gets the text code "ali" instead of a number.
# The bottom rows (35-37: fat/she/taj admins) already carry their own
# text codes and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..21 -> sequential numbers 1..20
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 3).Value = $row - 1
}

# Row 22 -> special text code instead of a number
$ws.Cells.Item(22, 3).Value = "ali"

# Rows 23..34 -> sequential numbers 22..33 (continuing the same series)
for ($row = 23; $row -le 34; $row++) {
    $ws.Cells.Item($row, 3).Value = $row - 1
}

# Rows 35..37 keep their existing text codes (fat/she/taj) - untouched.

# Restore the view: active cell / selection on C22, scrolled so row 19
# is the first visible row (matches the saved worksheet view state).
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
[void]$ws.Range("C22").Select()
